# Updated cryptos list values (Price / Volume(1h) columns) as scraped on
# Sat Oct  7 19:25:57 UTC 2023.
#
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text in
# the workbook (e.g. "27.895.91", "  -0.19%  "). Several of the new Price
# values look like ordinary decimals (e.g. "212.01"); left to its own devices
# Excel's Range.Value setter would helpfully reinterpret those as numbers and
# silently drop meaningful trailing zeros (e.g. "0.870" -> 0.87). Force those
# specific cells to Text format first so the literal string is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D8", "D9", "D11", "D14", "D15", "D16", "D18", "D23", "D25", "D26", "D27", "D30", "D35", "D39", "D40", "D42", "D43", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    "D2" = "27.895.91"
    "E2" = "  -0.19%  "
    "D3" = "1.633.88"
    "E3" = "  -0.91%  "
    "E4" = "  +0.01%  "
    "D5" = "212.01"
    "E6" = "  -0.78%  "
    "E7" = "  +0.02%  "
    "D8" = "23.16"
    "E8" = "  -1.69%  "
    "D9" = "0.257"
    "E9" = "  -3.46%  "
    "E10" = "  -0.40%  "
    "D11" = "0.0879"
    "E11" = "  +0.82%  "
    "D12" = "1.866.04"
    "E12" = "  -0.84%  "
    "D13" = "1.632.07"
    "E13" = "  -0.87%  "
    "D14" = "4.05"
    "E14" = "  -0.72%  "
    "D15" = "0.567"
    "E15" = "  +0.53%  "
    "D16" = "65.11"
    "E16" = "  -0.87%  "
    "D17" = "27.897.60"
    "E17" = "  -0.15%  "
    "D18" = "229.58"
    "E18" = "  -1.27%  "
    "D19" = "0.0₃0719"
    "E19" = "  -0.56%  "
    "E20" = "  -2.39%  "
    "E21" = "  -0.01%  "
    "E22" = "  -0.75%  "
    "D23" = "10.37"
    "E23" = "  -3.33%  "
    "E24" = "  -3.95%  "
    "D25" = "152.52"
    "E25" = "  +0.16%  "
    "D26" = "6.94"
    "E26" = "  +0.53%  "
    "D27" = "15.60"
    "E27" = "  -0.86%  "
    "E28" = "  -0.73%  "
    "E29" = "  +0.00%  "
    "D30" = "1.18"
    "E30" = "  -0.84%  "
    "E31" = "  -0.79%  "
    "E32" = "  +0.48%  "
    "D33" = "1.398.57"
    "E33" = "  -3.82%  "
    "E34" = "  -1.60%  "
    "D35" = "1.55"
    "E35" = "  +0.03%  "
    "E36" = "  +8.82%  "
    "E37" = "  +1.52%  "
    "E38" = "  +0.14%  "
    "D39" = "0.559"
    "E39" = "  -0.28%  "
    "D40" = "0.870"
    "E40" = "  -2.35%  "
    "E41" = "  -0.09%  "
    "D42" = "0.999"
    "D43" = "66.75"
    "E43" = "  -3.73%  "
    "E44" = "  +2.65%  "
    "E45" = "  +0.90%  "
    "E46" = "  -1.53%  "
    "D47" = "1.774.81"
    "E47" = "  -0.90%  "
    "D48" = "87.66"
    "E48" = "  -1.23%  "
    "E49" = "  -0.59%  "
    "E50" = "  -0.12%  "
    "D51" = "7.54"
    "E51" = "  -2.11%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

